# module col dynamic change
# Updates the "dynamic" module data cells (phone numbers / dates / record id)
# on each scenario sheet with freshly generated values.

$wb = $excel.ActiveWorkbook

# Sheet1 ("AddNew" scenario) - all dynamic columns including Date / Date&Time / User1RecId
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F2").Value = "5645901384"
$ws1.Range("N2").Value = "2024-02-16"
$ws1.Range("P2").Value = "2024-02-16 03:15:20 PM"
$ws1.Range("AC2").Value = "2024-02-16"
$ws1.Range("AE2").Value = "8534621837"
$ws1.Range("AN2").Value = "94434"
$ws1.Range("AT2").Value = "1407798439"
$ws1.Range("AX2").Value = "7842136796"

# Sheet2 ("SummaryAdd" scenario) - phone number columns only
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("F2").Value = "5645901384"
$ws2.Range("AE2").Value = "8534621837"
$ws2.Range("AT2").Value = "1407798439"
$ws2.Range("AX2").Value = "7842136796"

# Sheet3 ("Duplicate" scenario) - phone number columns only
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("F2").Value = "5645901384"
$ws3.Range("AE2").Value = "8534621837"
$ws3.Range("AT2").Value = "1407798439"
$ws3.Range("AX2").Value = "7842136796"

# Sheet4 ("EditSudhakar" scenario) - phone number columns only
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("F2").Value = "5645901384"
$ws4.Range("AE2").Value = "8534621837"
$ws4.Range("AT2").Value = "1407798439"
$ws4.Range("AX2").Value = "7842136796"
